$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 55
$ws.Cells.Item(55, 8).Value = 179.71428
$ws.Cells.Item(55, 9).Value = 136.66667
$ws.Cells.Item(55, 10).Value = 212
$ws.Cells.Item(55, 11).Value = 136.66667
$ws.Cells.Item(55, 12).Value = 212
$ws.Cells.Item(55, 13).Value = 77.33332999999999
$ws.Cells.Item(55, 14).Value = -640
# Row 129
$ws.Cells.Item(129, 8).Value = 916.2
$ws.Cells.Item(129, 10).Value = 958
$ws.Cells.Item(129, 12).Value = 2874
$ws.Cells.Item(129, 14).Value = -12874
# Row 137
$ws.Cells.Item(137, 8).Value = 1166.1915
$ws.Cells.Item(137, 9).Value = 864.45715
$ws.Cells.Item(137, 10).Value = 2046.25
$ws.Cells.Item(137, 11).Value = 2593.37145
$ws.Cells.Item(137, 12).Value = 6138.75
$ws.Cells.Item(137, 13).Value = -43.37144999999964
$ws.Cells.Item(137, 14).Value = -11238.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Cells.Item(45, 8).Value = 1575.375
$ws.Cells.Item(45, 9).Value = 1098
$ws.Cells.Item(45, 10).Value = 2371
$ws.Cells.Item(45, 11).Value = 1098
$ws.Cells.Item(45, 12).Value = 2371
$ws.Cells.Item(45, 13).Value = -721
$ws.Cells.Item(45, 14).Value = -3125
# Row 61
$ws.Cells.Item(61, 8).Value = 5209.091
$ws.Cells.Item(61, 9).Value = 3033.3333
$ws.Cells.Item(61, 11).Value = 3033.3333
$ws.Cells.Item(61, 13).Value = -2821.3333
# Row 74
$ws.Cells.Item(74, 8).Value = 1041.4103
$ws.Cells.Item(74, 9).Value = 1208.8077
$ws.Cells.Item(74, 10).Value = 706.61536
$ws.Cells.Item(74, 11).Value = 1208.8077
$ws.Cells.Item(74, 12).Value = 706.61536
$ws.Cells.Item(74, 13).Value = -334.8077000000001
$ws.Cells.Item(74, 14).Value = -2454.61536
# Row 77
$ws.Cells.Item(77, 8).Value = 1041.4103
$ws.Cells.Item(77, 9).Value = 1208.8077
$ws.Cells.Item(77, 10).Value = 706.61536
$ws.Cells.Item(77, 11).Value = 6044.038500000001
$ws.Cells.Item(77, 12).Value = 3533.0768
$ws.Cells.Item(77, 13).Value = -1676.038500000001
$ws.Cells.Item(77, 14).Value = -12269.0768
# Row 122
$ws.Cells.Item(122, 8).Value = 1617.75
$ws.Cells.Item(122, 9).Value = 1526.625
$ws.Cells.Item(122, 10).Value = 1800
$ws.Cells.Item(122, 11).Value = 4579.875
$ws.Cells.Item(122, 12).Value = 5400
$ws.Cells.Item(122, 13).Value = -2129.875
$ws.Cells.Item(122, 14).Value = -10300
# Row 123
$ws.Cells.Item(123, 8).Value = 24304.166
$ws.Cells.Item(123, 10).Value = 24304.166
$ws.Cells.Item(123, 12).Value = 24304.166
$ws.Cells.Item(123, 14).Value = -34104.166
# Row 131
$ws.Cells.Item(131, 8).Value = 0
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(131, 14).ClearContents()
# Row 136
$ws.Cells.Item(136, 8).Value = 5209.091
$ws.Cells.Item(136, 9).Value = 3033.3333
$ws.Cells.Item(136, 11).Value = 9099.999899999999
$ws.Cells.Item(136, 13).Value = -6549.999899999999
# Row 139
$ws.Cells.Item(139, 8).Value = 55786.35
$ws.Cells.Item(139, 10).Value = 55335.5
$ws.Cells.Item(139, 12).Value = 55335.5
$ws.Cells.Item(139, 14).Value = -65615.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 122
$ws.Cells.Item(122, 8).Value = 41800
$ws.Cells.Item(122, 10).Value = 41800
$ws.Cells.Item(122, 12).Value = 41800
$ws.Cells.Item(122, 14).Value = -51600
# Row 134
$ws.Cells.Item(134, 8).Value = 2353.7874
$ws.Cells.Item(134, 9).Value = 2054.4375
$ws.Cells.Item(134, 11).Value = 6163.3125
$ws.Cells.Item(134, 13).Value = -3628.3125

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Cells.Item(58, 8).Value = 1124276.1
$ws.Cells.Item(58, 9).Value = 1951139.6
$ws.Cells.Item(58, 10).Value = 2104.3572
$ws.Cells.Item(58, 11).Value = 1951139.6
$ws.Cells.Item(58, 12).Value = 2104.3572
$ws.Cells.Item(58, 13).Value = -1950936.6
$ws.Cells.Item(58, 14).Value = -2510.3572
# Row 134
$ws.Cells.Item(134, 8).Value = 1661.375
$ws.Cells.Item(134, 9).Value = 1403.5217
$ws.Cells.Item(134, 10).Value = 2320.3333
$ws.Cells.Item(134, 11).Value = 4210.5651
$ws.Cells.Item(134, 12).Value = 6960.999899999999
$ws.Cells.Item(134, 13).Value = -1675.5651
$ws.Cells.Item(134, 14).Value = -12030.9999
# Row 136
$ws.Cells.Item(136, 8).Value = 1124276.1
$ws.Cells.Item(136, 9).Value = 1951139.6
$ws.Cells.Item(136, 10).Value = 2104.3572
$ws.Cells.Item(136, 11).Value = 5853418.800000001
$ws.Cells.Item(136, 12).Value = 6313.071599999999
$ws.Cells.Item(136, 13).Value = -5850868.800000001
$ws.Cells.Item(136, 14).Value = -11413.0716

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Cells.Item(113, 8).Value = 557.2857
$ws.Cells.Item(113, 9).Value = 519.125
$ws.Cells.Item(113, 10).Value = 580.7692
$ws.Cells.Item(113, 11).Value = 1557.375
$ws.Cells.Item(113, 12).Value = 1742.3076
$ws.Cells.Item(113, 13).Value = 612.625
$ws.Cells.Item(113, 14).Value = -6082.3076
# Row 129
$ws.Cells.Item(129, 8).Value = 3334373
$ws.Cells.Item(129, 9).Value = 707.5
$ws.Cells.Item(129, 10).Value = 4546615
$ws.Cells.Item(129, 11).Value = 2122.5
$ws.Cells.Item(129, 12).Value = 13639845
$ws.Cells.Item(129, 13).Value = 2877.5
$ws.Cells.Item(129, 14).Value = -13649845
# Row 131
$ws.Cells.Item(131, 8).Value = 951.75
$ws.Cells.Item(131, 10).Value = 1009.0659
$ws.Cells.Item(131, 12).Value = 3027.1977
$ws.Cells.Item(131, 14).Value = -13107.1977

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 51
$ws.Cells.Item(51, 8).Value = 31372.77
$ws.Cells.Item(51, 10).Value = 31372.77
$ws.Cells.Item(51, 12).Value = 31372.77
$ws.Cells.Item(51, 14).Value = -32390.77
# Row 70
$ws.Cells.Item(70, 8).Value = 6238.276
$ws.Cells.Item(70, 9).Value = 5800.4375
$ws.Cells.Item(70, 10).Value = 6777.154
$ws.Cells.Item(70, 11).Value = 5800.4375
$ws.Cells.Item(70, 12).Value = 6777.154
$ws.Cells.Item(70, 13).Value = -5530.4375
$ws.Cells.Item(70, 14).Value = -7317.154
# Row 73
$ws.Cells.Item(73, 8).Value = 6238.276
$ws.Cells.Item(73, 9).Value = 5800.4375
$ws.Cells.Item(73, 10).Value = 6777.154
$ws.Cells.Item(73, 11).Value = 5800.4375
$ws.Cells.Item(73, 12).Value = 6777.154
$ws.Cells.Item(73, 13).Value = -4864.4375
$ws.Cells.Item(73, 14).Value = -8649.154
# Row 109
$ws.Cells.Item(109, 8).Value = 11211.066
$ws.Cells.Item(109, 10).Value = 11211.066
$ws.Cells.Item(109, 12).Value = 11211.066
$ws.Cells.Item(109, 14).Value = -13291.066
# Row 122
$ws.Cells.Item(122, 8).Value = 2841.2
$ws.Cells.Item(122, 9).Value = 2433.647
$ws.Cells.Item(122, 10).Value = 3707.25
$ws.Cells.Item(122, 11).Value = 7300.941
$ws.Cells.Item(122, 12).Value = 11121.75
$ws.Cells.Item(122, 13).Value = -4850.941
$ws.Cells.Item(122, 14).Value = -16021.75
# Row 123
$ws.Cells.Item(123, 8).Value = 9656.272000000001
$ws.Cells.Item(123, 10).Value = 9656.272000000001
$ws.Cells.Item(123, 12).Value = 9656.272000000001
$ws.Cells.Item(123, 14).Value = -14556.272

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 924
$ws.Cells.Item(22, 9).Value = 706.6667
$ws.Cells.Item(22, 11).Value = 706.6667
$ws.Cells.Item(22, 13).Value = -411.6667
# Row 27
$ws.Cells.Item(27, 8).Value = 924
$ws.Cells.Item(27, 9).Value = 706.6667
$ws.Cells.Item(27, 11).Value = 706.6667
$ws.Cells.Item(27, 13).Value = -599.6667
# Row 40
$ws.Cells.Item(40, 8).Value = 3848.7368
$ws.Cells.Item(40, 10).Value = 3274.1667
$ws.Cells.Item(40, 12).Value = 3274.1667
$ws.Cells.Item(40, 14).Value = -3546.1667
# Row 46
$ws.Cells.Item(46, 8).Value = 1100
$ws.Cells.Item(46, 9).Value = 666.6667
$ws.Cells.Item(46, 10).Value = 1750
$ws.Cells.Item(46, 11).Value = 666.6667
$ws.Cells.Item(46, 12).Value = 1750
$ws.Cells.Item(46, 13).Value = -478.6667
$ws.Cells.Item(46, 14).Value = -2126
# Row 122
$ws.Cells.Item(122, 8).Value = 14067797
$ws.Cells.Item(122, 9).Value = 15629657
$ws.Cells.Item(122, 10).Value = 12505937
$ws.Cells.Item(122, 11).Value = 46888971
$ws.Cells.Item(122, 12).Value = 37517811
$ws.Cells.Item(122, 13).Value = -46886521
$ws.Cells.Item(122, 14).Value = -37522711

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Cells.Item(62, 8).Value = 4312.875
$ws.Cells.Item(62, 9).Value = 3940
$ws.Cells.Item(62, 11).Value = 3940
$ws.Cells.Item(62, 13).Value = -3316
# Row 65
$ws.Cells.Item(65, 8).Value = 4312.875
$ws.Cells.Item(65, 9).Value = 3940
$ws.Cells.Item(65, 11).Value = 19700
$ws.Cells.Item(65, 13).Value = -16580
# Row 122
$ws.Cells.Item(122, 8).Value = 10778423
$ws.Cells.Item(122, 9).Value = 15627073
$ws.Cells.Item(122, 10).Value = 4810853.5
$ws.Cells.Item(122, 11).Value = 46881219
$ws.Cells.Item(122, 12).Value = 14432560.5
$ws.Cells.Item(122, 13).Value = -46878769
$ws.Cells.Item(122, 14).Value = -14437460.5
# Row 123
$ws.Cells.Item(123, 8).Value = 22272.154
$ws.Cells.Item(123, 10).Value = 22272.154
$ws.Cells.Item(123, 12).Value = 22272.154
$ws.Cells.Item(123, 14).Value = -32072.154
# Row 136
$ws.Cells.Item(136, 8).Value = 1650.9546
$ws.Cells.Item(136, 9).Value = 1759.0526
$ws.Cells.Item(136, 11).Value = 5277.1578
$ws.Cells.Item(136, 13).Value = -2727.1578

"done"